$d = $word.ActiveDocument
$sec = $d.Sections.First

# The Pearson Edexcel logo picture appears in both the default footer
# (index 1) and the first-page footer (index 2) of the document's only
# section; its inline-shape "Name" (OOXML <wp:docPr name="...">) flips
# from image1.png to image2.png in both places.
for ($fi = 1; $fi -le 3; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        for ($k = 1; $k -le $ftr.Range.InlineShapes.Count; $k++) {
            $ftr.Range.InlineShapes.Item($k).Name = "image2.png"
        }
    }
}

# The BTEC logo picture sits in the first-page header (index 2); its name
# flips the other way, from image2.jpg to image1.jpg.
for ($hi = 1; $hi -le 3; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        for ($k = 1; $k -le $hdr.Range.InlineShapes.Count; $k++) {
            $hdr.Range.InlineShapes.Item($k).Name = "image1.jpg"
        }
    }
}
